{"js": "// Fix a couple of label-name typos in the handout:\n//  1. \"Lab4Main\" -> \"Main\"\n//  2. \"look left\"  -> \"look robot's left\"\n//  3. \"look right\" -> \"look robot's right\"\n// (curly apostrophe U+2019, matching the author's text)\n\nconst body = context.document.body;\n\n// 1) Lab4Main -> Main\nconst mainResults = body.search(\"Lab4Main\", { matchCase: true });\nmainResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < mainResults.items.length; i++) {\n  mainResults.items[i].insertText(\"Main\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) look left -> look robot's left\nconst leftResults = body.search(\"look left\", { matchCase: true });\nleftResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < leftResults.items.length; i++) {\n  leftResults.items[i].insertText(\"look robot\\u2019s left\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) look right -> look robot's right\nconst rightResults = body.search(\"look right\", { matchCase: true });\nrightResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < rightResults.items.length; i++) {\n  rightResults.items[i].insertText(\"look robot\\u2019s right\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix a couple of label-name typos in the handout:\n#  1. \"Lab4Main\" -> \"Main\"\n#  2. \"look left\"  -> \"look robot's left\"   (curly apostrophe, U+2019)\n#  3. \"look right\" -> \"look robot's right\"  (curly apostrophe, U+2019)\n\n$d = $word.ActiveDocument\n$apos = [char]0x2019\n\n# 1) Lab4Main -> Main\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Lab4Main\"\n$find.Replacement.Text = \"Main\"\n$find.Execute([ref]$find.Text, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n\n# 2) look left -> look robot's left\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"look left\"\n$find2.Replacement.Text = \"look robot\" + $apos + \"s left\"\n$find2.Execute([ref]$find2.Text, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2)\n\n# 3) look right -> look robot's right\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"look right\"\n$find3.Replacement.Text = \"look robot\" + $apos + \"s right\"\n$find3.Execute([ref]$find3.Text, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$find3.Replacement.Text, 2)\n"}
